# Generate Report for Handoff
#
# A new handoff run produced a fresh GUID-named source snapshot and a new
# xliff content hash; refresh the localization-status report so every
# sheet (Overview, zh-cn, de-de) points at the new file names and the
# updated handoff/handback timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "9727361f-5031-460b-9a73-94786b7702fe"
$newGuid = "2c38df77-837e-43d9-b4c1-8276e797efe4"

$oldHash = "f0bfdaaf681f76cbc2760d5d8be1a950bd571b77"
$newHash = "820350ece49fe2ca1b5fc060f04ed6e44c428dec"

$newFileName = "$newGuid.md"
$newPathName = "e2e\$newGuid.md"
$newZhXlf    = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlf    = "$newGuid.$newHash.de-de.xlf"

$overviewDate = "2016-08-27 04:56:30"
$zhCnDate     = "2016-08-27 04:56:25"

# Hyperlink target URLs are unchanged by this edit - only the cached
# "display" text (and the cell text itself) needs to move to the new
# file name.
$overviewLinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/520fb0fd6a08e81499c112e77cff716ba92060cb/e2e/$oldGuid.md"
$zhCnLinkTarget      = $overviewLinkTarget
$deDeLinkTarget      = $overviewLinkTarget

function Update-Hyperlink($ws, $range, $target, $displayText) {
    # Dropping + re-adding is the only way this host lets us refresh the
    # cached display text of an existing hyperlink without leaving a
    # stale duplicate entry behind.
    $range.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($range, $target, "", "", $displayText) | Out-Null
}

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newPathName
Update-Hyperlink $wsOverview $wsOverview.Range("B2") $overviewLinkTarget $newPathName
$wsOverview.Range("G2").Value = $overviewDate

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newFileName
Update-Hyperlink $wsZhCn $wsZhCn.Range("A2") $zhCnLinkTarget $newFileName
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = $zhCnDate

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newFileName
Update-Hyperlink $wsDeDe $wsDeDe.Range("A2") $deDeLinkTarget $newFileName
$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = $overviewDate
